# Insert a new data row at row 224 (shifting existing rows 224-317 down to
# 225-318) and populate the new row with the latest "Femacal de La Calera -
# Berenjena" price record. This mirrors the author's weekly update, which
# prepends one fresh observation to the top of this sub-range.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Rows.Item(224).Insert()

$ws.Range("A224").Value = 3
$ws.Range("B224").Value = "Femacal de La Calera"
$ws.Range("C224").Value = "Coquimbo"
$ws.Range("D224").Value = 44755
$ws.Range("E224").Value = 5
$ws.Range("F224").Value = 100112001
$ws.Range("G224").Value = "Berenjena"
$ws.Range("H224").Value = "Sin especificar"
$ws.Range("I224").Value = "Primera"
$ws.Range("J224").Value = 75
$ws.Range("K224").Value = 8000
$ws.Range("L224").Value = 8000
$ws.Range("M224").Value = 8000
$ws.Range("N224").Value = "`$/caja 60 unidades"
$ws.Range("O224").Value = "Región de Arica y Parinacota"
$ws.Range("P224").Value = 133
$ws.Range("Q224").Value = 60
$ws.Range("R224").Value = "Hortaliza"
